# Auto-applies the "Cactuar_Profits" market-data refresh across all 8 job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR): updated Universalis price
# snapshots (currentAveragePrice / NQ / HQ, columns H-N) recompute per-row.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item(1)
# row 5
$ws.Cells.Item(5, 8).Value = 569
$ws.Cells.Item(5, 9).Value = 138.33333
$ws.Cells.Item(5, 11).Value = 138.33333
$ws.Cells.Item(5, 13).Value = -23.33332999999999
# row 43
$ws.Cells.Item(43, 8).Value = 1889.5454
$ws.Cells.Item(43, 9).Value = 1275
$ws.Cells.Item(43, 10).Value = 2026.1111
$ws.Cells.Item(43, 11).Value = 1275
$ws.Cells.Item(43, 12).Value = 2026.1111
$ws.Cells.Item(43, 13).Value = -1206
$ws.Cells.Item(43, 14).Value = -2164.1111
# row 53
$ws.Cells.Item(53, 8).Value = 291.88235
$ws.Cells.Item(53, 9).Value = 65.375
$ws.Cells.Item(53, 11).Value = 65.375
$ws.Cells.Item(53, 13).Value = 571.625
# row 80
$ws.Cells.Item(80, 8).Value = 37879320
$ws.Cells.Item(80, 10).Value = 83334350
$ws.Cells.Item(80, 12).Value = 250003050
$ws.Cells.Item(80, 14).Value = -250005046
# row 83
$ws.Cells.Item(83, 8).Value = 37879320
$ws.Cells.Item(83, 10).Value = 83334350
$ws.Cells.Item(83, 12).Value = 750009150
$ws.Cells.Item(83, 14).Value = -750019134
# row 100
$ws.Cells.Item(100, 8).Value = 2162.8462
$ws.Cells.Item(100, 9).Value = 1381.2
$ws.Cells.Item(100, 11).Value = 1381.2
$ws.Cells.Item(100, 13).Value = -840.2
# row 112
$ws.Cells.Item(112, 8).Value = 2968.2
$ws.Cells.Item(112, 10).Value = 3018.6924
$ws.Cells.Item(112, 12).Value = 9056.0772
$ws.Cells.Item(112, 14).Value = -11272.0772
# row 113
$ws.Cells.Item(113, 8).Value = 72462.3
$ws.Cells.Item(113, 9).Value = 2708.1667
$ws.Cells.Item(113, 11).Value = 2708.1667
$ws.Cells.Item(113, 13).Value = 545.8332999999998
# row 129
$ws.Cells.Item(129, 8).Value = 1464.3959
$ws.Cells.Item(129, 10).Value = 2344.4285
$ws.Cells.Item(129, 12).Value = 7033.2855
$ws.Cells.Item(129, 14).Value = -17033.2855
# row 132
$ws.Cells.Item(132, 8).Value = 282540.44
$ws.Cells.Item(132, 9).Value = 320171.22
$ws.Cells.Item(132, 10).Value = 19125
$ws.Cells.Item(132, 11).Value = 960513.6599999999
$ws.Cells.Item(132, 12).Value = 57375
$ws.Cells.Item(132, 13).Value = -957983.6599999999
$ws.Cells.Item(132, 14).Value = -62435
# row 135
$ws.Cells.Item(135, 8).Value = 2855.2
$ws.Cells.Item(135, 9).Value = 1614.2858
$ws.Cells.Item(135, 11).Value = 14528.5722
$ws.Cells.Item(135, 13).Value = -11993.5722
# row 138
$ws.Cells.Item(138, 8).Value = 2415.32
$ws.Cells.Item(138, 10).Value = 2532.4353
$ws.Cells.Item(138, 12).Value = 7597.3059
$ws.Cells.Item(138, 14).Value = -17877.3059
# --- ARM ---
$ws = $wb.Worksheets.Item(2)
# row 32
$ws.Cells.Item(32, 8).Value = 2009.91
$ws.Cells.Item(32, 9).Value = 2009.91
$ws.Cells.Item(32, 11).Value = 2009.91
$ws.Cells.Item(32, 13).Value = -1722.91
# row 61
$ws.Cells.Item(61, 8).Value = 6123.148
$ws.Cells.Item(61, 9).Value = 5705.0454
$ws.Cells.Item(61, 11).Value = 5705.0454
$ws.Cells.Item(61, 13).Value = -5493.0454
# row 102
$ws.Cells.Item(102, 8).Value = 2283.7693
$ws.Cells.Item(102, 9).Value = 2240.75
$ws.Cells.Item(102, 11).Value = 2240.75
$ws.Cells.Item(102, 13).Value = -618.75
# row 110
$ws.Cells.Item(110, 8).Value = 1403.6757
$ws.Cells.Item(110, 9).Value = 1182.8667
$ws.Cells.Item(110, 11).Value = 1182.8667
$ws.Cells.Item(110, 13).Value = 862.1333
# row 122
$ws.Cells.Item(122, 8).Value = 2222.0488
$ws.Cells.Item(122, 9).Value = 1328.0588
$ws.Cells.Item(122, 10).Value = 6564.2856
$ws.Cells.Item(122, 11).Value = 3984.1764
$ws.Cells.Item(122, 12).Value = 19692.8568
$ws.Cells.Item(122, 13).Value = -1534.1764
$ws.Cells.Item(122, 14).Value = -24592.8568
# row 130
$ws.Cells.Item(130, 8).Value = 12481.5
$ws.Cells.Item(130, 10).Value = 12481.5
$ws.Cells.Item(130, 12).Value = 12481.5
$ws.Cells.Item(130, 14).Value = -22521.5
# row 132
$ws.Cells.Item(132, 8).Value = 11408.843
$ws.Cells.Item(132, 9).Value = 11331.534
$ws.Cells.Item(132, 11).Value = 33994.602
$ws.Cells.Item(132, 13).Value = -31464.602
# row 136
$ws.Cells.Item(136, 8).Value = 6123.148
$ws.Cells.Item(136, 9).Value = 5705.0454
$ws.Cells.Item(136, 11).Value = 17115.1362
$ws.Cells.Item(136, 13).Value = -14565.1362
# --- BSM ---
$ws = $wb.Worksheets.Item(3)
# row 94
$ws.Cells.Item(94, 8).Value = 1577.28
$ws.Cells.Item(94, 9).Value = 1723.4375
$ws.Cells.Item(94, 11).Value = 1723.4375
$ws.Cells.Item(94, 13).Value = -1272.4375
# row 99
$ws.Cells.Item(99, 8).Value = 3504.2666
$ws.Cells.Item(99, 9).Value = 3742.182
$ws.Cells.Item(99, 11).Value = 3742.182
$ws.Cells.Item(99, 13).Value = -2244.182
# row 105
$ws.Cells.Item(105, 8).Value = 3004.1365
$ws.Cells.Item(105, 9).Value = 3034
$ws.Cells.Item(105, 11).Value = 3034
$ws.Cells.Item(105, 13).Value = -1287
# row 128
$ws.Cells.Item(128, 8).Value = 7583.75
$ws.Cells.Item(128, 9).Value = 7583.75
$ws.Cells.Item(128, 11).Value = 22751.25
$ws.Cells.Item(128, 13).Value = -20261.25
# row 134
$ws.Cells.Item(134, 8).Value = 3632.7727
$ws.Cells.Item(134, 9).Value = 3632.7727
$ws.Cells.Item(134, 11).Value = 10898.3181
$ws.Cells.Item(134, 13).Value = -8363.3181
# --- CRP ---
$ws = $wb.Worksheets.Item(4)
# row 16
$ws.Cells.Item(16, 8).Value = 1615.5
$ws.Cells.Item(16, 9).Value = 1539.1
$ws.Cells.Item(16, 10).Value = 1997.5
$ws.Cells.Item(16, 11).Value = 1539.1
$ws.Cells.Item(16, 12).Value = 1997.5
$ws.Cells.Item(16, 13).Value = -1252.1
$ws.Cells.Item(16, 14).Value = -2571.5
# row 93
$ws.Cells.Item(93, 8).Value = 55575676
$ws.Cells.Item(93, 9).Value = 7365.8335
$ws.Cells.Item(93, 10).Value = 166712300
$ws.Cells.Item(93, 11).Value = 7365.8335
$ws.Cells.Item(93, 12).Value = 166712300
$ws.Cells.Item(93, 13).Value = -5493.8335
$ws.Cells.Item(93, 14).Value = -166716044
# row 99
$ws.Cells.Item(99, 8).Value = 7080.0938
$ws.Cells.Item(99, 9).Value = 7745.75
$ws.Cells.Item(99, 10).Value = 5970.6665
$ws.Cells.Item(99, 11).Value = 7745.75
$ws.Cells.Item(99, 12).Value = 5970.6665
$ws.Cells.Item(99, 13).Value = -6247.75
$ws.Cells.Item(99, 14).Value = -8966.666499999999
# row 113
$ws.Cells.Item(113, 8).Value = 1615.5
$ws.Cells.Item(113, 9).Value = 1539.1
$ws.Cells.Item(113, 10).Value = 1997.5
$ws.Cells.Item(113, 11).Value = 1539.1
$ws.Cells.Item(113, 12).Value = 1997.5
$ws.Cells.Item(113, 13).Value = 630.9000000000001
$ws.Cells.Item(113, 14).Value = -6337.5
# row 114
$ws.Cells.Item(114, 8).Value = 53488
$ws.Cells.Item(114, 10).Value = 53488
$ws.Cells.Item(114, 12).Value = 53488
$ws.Cells.Item(114, 14).Value = -62166
# row 115
$ws.Cells.Item(115, 8).Value = 49999.266
$ws.Cells.Item(115, 10).Value = 49999.266
$ws.Cells.Item(115, 12).Value = 49999.266
$ws.Cells.Item(115, 14).Value = -52349.266
# row 126
$ws.Cells.Item(126, 8).Value = 7080.0938
$ws.Cells.Item(126, 9).Value = 7745.75
$ws.Cells.Item(126, 10).Value = 5970.6665
$ws.Cells.Item(126, 11).Value = 23237.25
$ws.Cells.Item(126, 12).Value = 17911.9995
$ws.Cells.Item(126, 13).Value = -20767.25
$ws.Cells.Item(126, 14).Value = -22851.9995
# row 134
$ws.Cells.Item(134, 8).Value = 3147.9375
$ws.Cells.Item(134, 9).Value = 3312
$ws.Cells.Item(134, 10).Value = 1999.5
$ws.Cells.Item(134, 11).Value = 9936
$ws.Cells.Item(134, 12).Value = 5998.5
$ws.Cells.Item(134, 13).Value = -7401
$ws.Cells.Item(134, 14).Value = -11068.5
# --- CUL ---
$ws = $wb.Worksheets.Item(5)
# row 14
$ws.Cells.Item(14, 8).Value = 316.66666
$ws.Cells.Item(14, 9).Value = 316.66666
$ws.Cells.Item(14, 11).Value = 949.9999799999999
$ws.Cells.Item(14, 13).Value = -776.9999799999999
# row 59
$ws.Cells.Item(59, 8).Value = 5666.6665
$ws.Cells.Item(59, 10).Value = 5666.6665
$ws.Cells.Item(59, 12).Value = 16999.9995
$ws.Cells.Item(59, 14).Value = -18079.9995
# row 88
$ws.Cells.Item(88, 8).Value = 9200
$ws.Cells.Item(88, 10).Value = 9200
$ws.Cells.Item(88, 12).Value = 27600
$ws.Cells.Item(88, 14).Value = -28456
# row 91
$ws.Cells.Item(91, 8).Value = 9200
$ws.Cells.Item(91, 10).Value = 9200
$ws.Cells.Item(91, 12).Value = 27600
$ws.Cells.Item(91, 14).Value = -30564
# row 98
$ws.Cells.Item(98, 8).Value = 1160.8667
$ws.Cells.Item(98, 9).Value = 1560.75
$ws.Cells.Item(98, 10).Value = 1015.4545
$ws.Cells.Item(98, 11).Value = 4682.25
$ws.Cells.Item(98, 12).Value = 3046.3635
$ws.Cells.Item(98, 13).Value = -3184.25
$ws.Cells.Item(98, 14).Value = -6042.3635
# row 113
$ws.Cells.Item(113, 8).Value = 721.5
$ws.Cells.Item(113, 9).Value = 769.2
$ws.Cells.Item(113, 11).Value = 2307.6
$ws.Cells.Item(113, 13).Value = -137.6000000000004
# row 122
$ws.Cells.Item(122, 8).Value = 2030.2
$ws.Cells.Item(122, 9).Value = 1627.8334
$ws.Cells.Item(122, 10).Value = 2298.4443
$ws.Cells.Item(122, 11).Value = 14650.5006
$ws.Cells.Item(122, 12).Value = 20685.9987
$ws.Cells.Item(122, 13).Value = -12200.5006
$ws.Cells.Item(122, 14).Value = -25585.9987
# row 137
$ws.Cells.Item(137, 8).Value = 53573828
$ws.Cells.Item(137, 10).Value = 3999.3333
$ws.Cells.Item(137, 12).Value = 11997.9999
$ws.Cells.Item(137, 14).Value = -22197.9999
# --- GSM ---
$ws = $wb.Worksheets.Item(6)
# row 80
$ws.Cells.Item(80, 8).Value = 56073.316
$ws.Cells.Item(80, 9).Value = 103309.6
$ws.Cells.Item(80, 11).Value = 103309.6
$ws.Cells.Item(80, 13).Value = -102311.6
# row 83
$ws.Cells.Item(83, 8).Value = 56073.316
$ws.Cells.Item(83, 9).Value = 103309.6
$ws.Cells.Item(83, 11).Value = 516548
$ws.Cells.Item(83, 13).Value = -511556
# row 97
$ws.Cells.Item(97, 8).Value = 1791.4706
$ws.Cells.Item(97, 9).Value = 1389.9286
$ws.Cells.Item(97, 11).Value = 1389.9286
$ws.Cells.Item(97, 13).Value = -893.9286
# row 102
$ws.Cells.Item(102, 8).Value = 10204726
$ws.Cells.Item(102, 9).Value = 11338107
$ws.Cells.Item(102, 11).Value = 11338107
$ws.Cells.Item(102, 13).Value = -11336485
# row 122
$ws.Cells.Item(122, 8).Value = 253960.47
$ws.Cells.Item(122, 9).Value = 456588.3
$ws.Cells.Item(122, 10).Value = 6304.222
$ws.Cells.Item(122, 11).Value = 1369764.9
$ws.Cells.Item(122, 12).Value = 18912.666
$ws.Cells.Item(122, 13).Value = -1367314.9
$ws.Cells.Item(122, 14).Value = -23812.666
# row 123
$ws.Cells.Item(123, 8).Value = 41036.223
$ws.Cells.Item(123, 10).Value = 41036.223
$ws.Cells.Item(123, 12).Value = 41036.223
$ws.Cells.Item(123, 14).Value = -45936.223
# row 132
$ws.Cells.Item(132, 8).Value = 113744.22
$ws.Cells.Item(132, 9).Value = 144599.78
$ws.Cells.Item(132, 11).Value = 433799.34
$ws.Cells.Item(132, 13).Value = -431269.34
# --- LTW ---
$ws = $wb.Worksheets.Item(7)
# row 46
$ws.Cells.Item(46, 8).Value = 4465.4585
$ws.Cells.Item(46, 9).Value = 2223.3333
$ws.Cells.Item(46, 11).Value = 2223.3333
$ws.Cells.Item(46, 13).Value = -2035.3333
# row 61
$ws.Cells.Item(61, 8).Value = 2172.963
$ws.Cells.Item(61, 9).Value = 2122.077
$ws.Cells.Item(61, 10).Value = 3496
$ws.Cells.Item(61, 11).Value = 2122.077
$ws.Cells.Item(61, 12).Value = 3496
$ws.Cells.Item(61, 13).Value = -1920.077
$ws.Cells.Item(61, 14).Value = -3900
# row 100
$ws.Cells.Item(100, 8).Value = 4532.5
$ws.Cells.Item(100, 9).Value = 3755.3333
$ws.Cells.Item(100, 11).Value = 3755.3333
$ws.Cells.Item(100, 13).Value = -3214.3333
# row 113
$ws.Cells.Item(113, 8).Value = 2172.963
$ws.Cells.Item(113, 9).Value = 2122.077
$ws.Cells.Item(113, 10).Value = 3496
$ws.Cells.Item(113, 11).Value = 2122.077
$ws.Cells.Item(113, 12).Value = 3496
$ws.Cells.Item(113, 13).Value = 47.92299999999977
$ws.Cells.Item(113, 14).Value = -7836
# --- WVR ---
$ws = $wb.Worksheets.Item(8)
# row 136
$ws.Cells.Item(136, 8).Value = 5130.913
$ws.Cells.Item(136, 9).Value = 3500.3103
$ws.Cells.Item(136, 10).Value = 7912.5293
$ws.Cells.Item(136, 11).Value = 10500.9309
$ws.Cells.Item(136, 12).Value = 23737.5879
$ws.Cells.Item(136, 13).Value = -7950.930899999999
$ws.Cells.Item(136, 14).Value = -28837.5879
# row 140
$ws.Cells.Item(140, 8).Value = 99999.5
$ws.Cells.Item(140, 10).Value = 99999.5
$ws.Cells.Item(140, 12).Value = 99999.5
$ws.Cells.Item(140, 14).Value = -110359.5
